# Sample Project / Main.xlsx - "Rules" sheet edit.
#
# The decision table's last rule row (row 11) had its Rule-name cell (B11)
# change from the text "R40" to the text "1". The cell keeps storing a
# *string* (not a number) both before and after, so we force text entry
# with a leading apostrophe - exactly what a user typing into that cell in
# Excel would do to keep "1" from being auto-converted to a numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("B11").Value = "'1"
